$wb = $excel.ActiveWorkbook

# xlEdgeLeft = 7, xlEdgeTop = 8, xlEdgeBottom = 9, xlEdgeRight = 10
# xlContinuous = 1, xlLineStyleNone = -4142, xlThin = 2 (weight)

function Set-TopBottomBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
    $rng.Borders.Item(10).LineStyle = -4142
}

function Set-TopBottomRightBorder($rng) {
    $rng.ClearFormats()
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
}

# ---- Sheet "quality_comparison" ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# ---- Sheet "computational_comparison" ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
